$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-5 (value column B, and reordered labels in column A)
$ws.Range("B2").Value = 74

$ws.Range("A3").Value = "Correct"
$ws.Range("B3").Value = 35

$ws.Range("A4").Value = "Wrong_Entity_Event_as_NonEvent"
$ws.Range("B4").Value = 27

$ws.Range("B5").Value = 3

# Add new rows 6 and 7, copying the style of row 5's cells so formatting matches
$ws.Range("A5:B5").Copy() | Out-Null
$ws.Range("A6:B6").PasteSpecial(-4122) | Out-Null

$ws.Range("A6").Value = "Wrong_Tag_B_as_I"
$ws.Range("B6").Value = 1

$ws.Range("A6:B6").Copy() | Out-Null
$ws.Range("A7:B7").PasteSpecial(-4122) | Out-Null

$ws.Range("A7").Value = "Wrong_Tag_S_as_B"
$ws.Range("B7").Value = 1

$excel.CutCopyMode = 0
